# Update countries & provincias Spain
# Applies the edits described by the paises.xlsx diff:
#  - bumps the "Datos actualizados" timestamp
#  - updates a handful of totals (Estados Unidos, Reino Unido, Malta)
#  - re-sorts "Principado de Andorra" to sit right after "Bulgaria"
#    (pushing Bosnia y Herzegovina / Eslovaquia down one row each)
#  - re-sorts "Azerbaiyan" to sit right after "Albania"
#    (pushing Burkina Faso / Vietnam / Reunion down one row each)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 16:20"

# --- Simple numeric corrections -----------------------------------------
# Estados Unidos (row 4)
$ws.Range("E4").Value = 118359
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 2231

# Reino Unido (row 11)
$ws.Range("E11").Value = 18152
$ws.Range("G11").Value = 216
$ws.Range("H11").Value = 1235

# Malta (row 96)
$ws.Range("F96").Value = 4

# --- Reorder: Principado de Andorra now right after Bulgaria (row 73) ---
# Row 74: now Andorra, with fresh totals
$ws.Range("A74").Value = "Principado de Andorra"
$ws.Range("B74").Value = 334
$ws.Range("C74").Value = 26
$ws.Range("D74").Value = 1
$ws.Range("E74").Value = 327
$ws.Range("F74").Value = 10
$ws.Range("G74").Value = 3
$ws.Range("H74").Value = 6

# Row 75: now Bosnia y Herzegovina (shifted down from old row 74)
$ws.Range("A75").Value = "Bosnia y Herzegovina"
$ws.Range("B75").Value = 323
$ws.Range("C75").Value = 45
$ws.Range("D75").Value = 8
$ws.Range("E75").Value = 309
$ws.Range("F75").Value = 1
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 6

# Row 76: now Eslovaquia (shifted down from old row 75)
$ws.Range("A76").Value = "Eslovaquia"
$ws.Range("B76").Value = 314
$ws.Range("C76").Value = 22
$ws.Range("D76").Value = 2
$ws.Range("E76").Value = 312
$ws.Range("F76").Value = 1
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 0

# Row 77 (Uruguay) is unchanged.

# --- Reorder: Azerbaiyan now right after Albania (row 87) ---------------
# Row 88: now Azerbaiyan, with fresh totals
$ws.Range("A88").Value = "Azerbaiyan"
$ws.Range("B88").Value = 209
$ws.Range("C88").Value = 27
$ws.Range("D88").Value = 15
$ws.Range("E88").Value = 190
$ws.Range("F88").Value = 23
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 4

# Row 89: now Burkina Faso (shifted down from old row 88)
$ws.Range("A89").Value = "Burkina Faso"
$ws.Range("B89").Value = 207
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 21
$ws.Range("E89").Value = 175
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 11

# Row 90: now Vietnam (shifted down from old row 89)
$ws.Range("A90").Value = "Vietnam"
$ws.Range("B90").Value = 188
$ws.Range("C90").Value = 14
$ws.Range("D90").Value = 21
$ws.Range("E90").Value = 167
$ws.Range("F90").Value = 3
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 0

# Row 91: now Reunion (shifted down from old row 90)
$ws.Range("A91").Value = "Reunion"
$ws.Range("B91").Value = 183
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 1
$ws.Range("E91").Value = 182
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 0

# Row 92 (Republica de Chipre) is unchanged.
